# Updated run for publication
# Applies refreshed frequency-table values (rows 2-5, columns B-X) to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{}
$updates[2] = @{
    2 = 0.00105485232067511   # B2
    3 = 0.992465340566606   # C2
    4 = 0.00406871609403255   # D2
    6 = 0.00120554550934298   # F2
    7 = 0.00406871609403255   # G2
    8 = 0.0764014466546112   # H2
    9 = 0.00979505726341169   # I2
    10 = 0.0748945147679325   # J2
    11 = 0.139240506329114   # K2
    12 = 0.000602772754671489   # L2
    13 = 0.997438215792646   # M2
    15 = 0.956298975286317   # O2
    16 = 0.000452079566003617   # P2
    18 = 0.000452079566003617   # R2
    19 = 0.00376732971669681   # S2
    20 = 0.996986136226643   # T2
    21 = 0.000150693188667872   # U2
    22 = 0.000150693188667872   # V2
    23 = 0.0322483423749247   # W2
    24 = 0.00165762507534659   # X2
}
$updates[3] = @{
    3 = 0.000602772754671489   # C3
    4 = 0.00105485232067511   # D3
    5 = 0.00467148884870404   # E3
    6 = 0.00376732971669681   # F3
    8 = 0.00105485232067511   # H3
    11 = 0.0229053646775166   # K3
    12 = 0.000904159132007233   # L3
    14 = 0.960066305003014   # N3
    16 = 0.0110006027727547   # P3
    17 = 0.000602772754671489   # Q3
    18 = 0.000452079566003617   # R3
    19 = 0.995328511151296   # S3
    20 = 0.000150693188667872   # T3
    23 = 0.00376732971669681   # W3
    24 = 0.00467148884870404   # X3
}
$updates[4] = @{
    2 = 0.993670886075949   # B4
    3 = 0.00180831826401447   # C4
    5 = 0.000301386377335744   # E4
    6 = 0.994575045207957   # F4
    7 = 0.9957805907173   # G4
    8 = 0.917420132610006   # H4
    9 = 0.985533453887884   # I4
    10 = 0.922845087402049   # J4
    11 = 0.794153104279687   # K4
    12 = 0.000452079566003617   # L4
    13 = 0.00256178420735383   # M4
    15 = 0.000301386377335744   # O4
    16 = 0.988245931283906   # P4
    17 = 0.00105485232067511   # Q4
    18 = 0.998794454490657   # R4
    20 = 0.00210970464135021   # T4
    21 = 0.999698613622664   # U4
    22 = 0.999849306811332   # V4
    23 = 0.95750452079566   # W4
    24 = 0.991862567811935   # X4
}
$updates[5] = @{
    2 = 0.00482218203737191   # B5
    3 = 0.00482218203737191   # C5
    4 = 0.994725738396624   # D5
    5 = 0.993972272453285   # E5
    6 = 0.000452079566003617   # F5
    8 = 0.000452079566003617   # H5
    9 = 0.00376732971669681   # I5
    11 = 0.0272754671488849   # K5
    12 = 0.99789029535865   # L5
    14 = 0.0375226039783002   # N5
    15 = 0.0429475587703436   # O5
    17 = 0.997588908981314   # Q5
    18 = 0.000150693188667872   # R5
    19 = 0.000602772754671489   # S5
    20 = 0.000602772754671489   # T5
    23 = 0.000150693188667872   # W5
    24 = 0.000301386377335744   # X5
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col]
    }
}

Write-Output "Updated $($updates.Count) rows of frequency table values"
